$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.144.81"
$ws.Range("E2").Value = "  +5.80%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.919.40"
$ws.Range("E3").Value = "  +2.60%  "
$ws.Range("E4").Value = "  -0.71%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "329.95"
$ws.Range("E5").Value = "  +4.57%  "
$ws.Range("E6").Value = "  -0.61%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5218"
$ws.Range("E7").Value = "  +2.77%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4091"
$ws.Range("E8").Value = "  +4.87%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08554"
$ws.Range("E9").Value = "  +2.41%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "43.05"
$ws.Range("E10").Value = "  +2.84%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.42"
$ws.Range("E12").Value = "  +9.83%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.415"
$ws.Range("E13").Value = "  +3.25%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.924.50"
$ws.Range("E14").Value = "  +2.99%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.415"
$ws.Range("E15").Value = "  +2.06%  "
$ws.Range("E16").Value = "  -0.78%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "95.45"
$ws.Range("E17").Value = "  +4.72%  "
$ws.Range("E18").Value = "  +1.38%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06689"
$ws.Range("E19").Value = "  -0.71%  "
$ws.Range("E20").Value = "  +4.11%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.000"
$ws.Range("E21").Value = "  -0.60%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.013"
$ws.Range("E22").Value = "  +1.58%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "30.144.49"
$ws.Range("E23").Value = "  +5.67%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.41"
$ws.Range("E24").Value = "  +2.95%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.208"
$ws.Range("E25").Value = "  +0.89%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.140.76"
$ws.Range("E26").Value = "  +2.70%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.11"
$ws.Range("E27").Value = "  +2.41%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "159.91"
$ws.Range("E28").Value = "  +0.86%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.446"
$ws.Range("E29").Value = "  +0.78%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "129.17"
$ws.Range("E30").Value = "  +2.31%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.086"
$ws.Range("E31").Value = "  +3.87%  "
$ws.Range("E32").Value = "  +2.28%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.062"
$ws.Range("E33").Value = "  +5.82%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.636"
$ws.Range("E34").Value = "  +0.64%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02495"
$ws.Range("E35").Value = "  +1.48%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06610"
$ws.Range("E36").Value = "  +0.30%  "
$ws.Range("E37").Value = "  +1.98%  "
$ws.Range("E38").Value = "  +4.55%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.184"
$ws.Range("E39").Value = "  +2.97%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.908"
$ws.Range("E40").Value = "  +0.17%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6560"
$ws.Range("E41").Value = "  +3.00%  "
$ws.Range("E42").Value = "  +1.16%  "
$ws.Range("E43").Value = "  +4.90%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6159"
$ws.Range("E44").Value = "  +2.60%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.20"
$ws.Range("E45").Value = "  +0.63%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.762"
$ws.Range("E46").Value = "  +2.25%  "
$ws.Range("E47").Value = "  +3.77%  "
$ws.Range("E48").Value = "  +3.00%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "124.60"
$ws.Range("E49").Value = "  +1.76%  "
$ws.Range("E50").Value = "  +10.47%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "79.90"
$ws.Range("E51").Value = "  +4.45%  "
